$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and volume-change (E) figures.
# Column D cells are forced to Text format first so numeric-looking
# price strings (e.g. "216.40") keep their exact original formatting
# instead of being auto-converted to a number by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.037.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.96"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.40"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5093"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2649"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06386"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.78"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07447"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.666.27"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.502"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5820"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008530"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.105.49"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.77"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.13"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.184"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.38"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.596"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1200"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06637"
$ws.Range("E28").Value = "  +14.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.335"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.316"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.546"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.508"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.656"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.017"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6120"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.369"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.307"
$ws.Range("E38").Value = "  +7.57%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01594"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8721"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.00"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.815.85"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("E45").Value = "  -4.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.35"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.009"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.064"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05227"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4285"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.025"
$ws.Range("E51").Value = "  +3.08%  "
